$d = $word.ActiveDocument

# Small line-number updates within the stack trace text (JavaMethodService line, etc.)
$old1 = "JavaMethodService.java:163)"
$new1 = "JavaMethodService.java:162)"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) { Write-Host "NOT FOUND: $old1" }

$old2 = "AbstractService.java:136)"
$new2 = "AbstractService.java:135)"
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) { Write-Host "NOT FOUND: $old2" }

$old3 = "EvaluationServices.java:168)"
$new3 = "EvaluationServices.java:172)"
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found3) { Write-Host "NOT FOUND: $old3" }

$old4 = "AstEvaluator.java:189)"
$new4 = "AstEvaluator.java:186)"
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $found4) { Write-Host "NOT FOUND: $old4" }

$old5 = "AstSwitch.java:118)"
$new5 = "AstSwitch.java:119)"
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $found5) { Write-Host "NOT FOUND: $old5" }

$old6 = "AstEvaluator.java:333)"
$new6 = "AstEvaluator.java:342)"
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
if (-not $found6) { Write-Host "NOT FOUND: $old6" }

$old7 = "AstSwitch.java:259)"
$new7 = "AstSwitch.java:262)"
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
if (-not $found7) { Write-Host "NOT FOUND: $old7" }

$old8 = "AstEvaluator.java:112)"
$new8 = "AstEvaluator.java:109)"
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
if (-not $found8) { Write-Host "NOT FOUND: $old8" }

$old9 = "GeneratedMethodAccessor74"
$new9 = "GeneratedMethodAccessor73"
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
if (-not $found9) { Write-Host "NOT FOUND: $old9" }

# Replace the large tail of the stack trace (surefire/equinox/tycho frames)
# with the Eclipse JDT JUnit runner frames.
$oldBlock = "	at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)
	at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)
	at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)
	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)
	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)
	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)
	at java.lang.reflect.Method.invoke(Method.java:498)
	at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)
	at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)
	at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)
	at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)
	at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)
	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)
	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)
	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)
	at java.lang.reflect.Method.invoke(Method.java:498)
	at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)
	at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)
	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)
	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)
	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)
	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)
	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)
	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)
	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)
	at java.lang.reflect.Method.invoke(Method.java:498)
	at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)
	at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)
	at org.eclipse.equinox.launcher.Main.run(Main.java:1498)
	at org.eclipse.equinox.launcher.Main.main(Main.java:1471)"
$newBlock = "	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)
	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$foundBlock = $d.Content.Find.Execute($oldBlock, $true, $false, $false, $false, $false, $true, 1, $false, $newBlock, 2)
if (-not $foundBlock) { Write-Host "NOT FOUND: oldBlock" }

Write-Host "Done"
